$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2381.5833
$ws.Range("J49").Value = 3247.375
$ws.Range("L49").Value = 9742.125
$ws.Range("N49").Value = -10014.125

$ws.Range("H92").Value = 266.08334
$ws.Range("I92").Value = 299.3
$ws.Range("K92").Value = 299.3
$ws.Range("M92").Value = 948.7

$ws.Range("H96").Value = 455.5
$ws.Range("I96").Value = 337.2143
$ws.Range("K96").Value = 1011.6429
$ws.Range("M96").Value = 361.3571000000001

$ws.Range("H129").Value = 817.88
$ws.Range("I129").Value = 445.5
$ws.Range("J129").Value = 993.1177
$ws.Range("K129").Value = 1336.5
$ws.Range("L129").Value = 2979.3531
$ws.Range("M129").Value = 3663.5
$ws.Range("N129").Value = -12979.3531

$ws.Range("H137").Value = 671913.2
$ws.Range("I137").Value = 1026806.94
$ws.Range("J137").Value = 6487.5
$ws.Range("K137").Value = 3080420.82
$ws.Range("L137").Value = 19462.5
$ws.Range("M137").Value = -3077870.82
$ws.Range("N137").Value = -24562.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22450.707
$ws.Range("I32").Value = 6422.9316
$ws.Range("J32").Value = 152453.78
$ws.Range("K32").Value = 6422.9316
$ws.Range("L32").Value = 152453.78
$ws.Range("M32").Value = -6135.9316
$ws.Range("N32").Value = -153027.78

$ws.Range("H61").Value = 1136.5
$ws.Range("I61").Value = 1080.7273
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 1080.7273
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -868.7273
$ws.Range("N61").Value = -2174

$ws.Range("H62").Value = 8226
$ws.Range("I62").Value = 8226
$ws.Range("K62").Value = 8226
$ws.Range("M62").Value = -7602

$ws.Range("H65").Value = 8226
$ws.Range("I65").Value = 8226
$ws.Range("K65").Value = 24678
$ws.Range("M65").Value = -21558

$ws.Range("H97").Value = 1500
$ws.Range("I97").Value = 780
$ws.Range("K97").Value = 780
$ws.Range("M97").Value = -284

$ws.Range("H136").Value = 1136.5
$ws.Range("I136").Value = 1080.7273
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 3242.1819
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = -692.1819
$ws.Range("N136").Value = -10350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23247.562
$ws.Range("J82").Value = 30871
$ws.Range("L82").Value = 30871
$ws.Range("N82").Value = -31637

$ws.Range("H85").Value = 23247.562
$ws.Range("J85").Value = 30871
$ws.Range("L85").Value = 30871
$ws.Range("N85").Value = -33523

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3345.037
$ws.Range("I31").Value = 1825.75
$ws.Range("J31").Value = 3984.7368
$ws.Range("K31").Value = 1825.75
$ws.Range("L31").Value = 3984.7368
$ws.Range("M31").Value = -1530.75
$ws.Range("N31").Value = -4574.736800000001

$ws.Range("H34").Value = 3345.037
$ws.Range("I34").Value = 1825.75
$ws.Range("J34").Value = 3984.7368
$ws.Range("K34").Value = 1825.75
$ws.Range("L34").Value = 3984.7368
$ws.Range("M34").Value = -1623.75
$ws.Range("N34").Value = -4388.736800000001

$ws.Range("H41").Value = 16815
$ws.Range("J41").Value = 20941
$ws.Range("L41").Value = 20941
$ws.Range("N41").Value = -21797

$ws.Range("H51").Value = 9102.6
$ws.Range("J51").Value = 9102.6
$ws.Range("L51").Value = 9102.6
$ws.Range("N51").Value = -10574.6

$ws.Range("H60").Value = 23297.066
$ws.Range("J60").Value = 24818.285
$ws.Range("L60").Value = 24818.285
$ws.Range("N60").Value = -25840.285

$ws.Range("H61").Value = 9102.6
$ws.Range("J61").Value = 9102.6
$ws.Range("L61").Value = 9102.6
$ws.Range("N61").Value = -9798.6

$ws.Range("H68").Value = 17398.666
$ws.Range("J68").Value = 17398.666
$ws.Range("L68").Value = 17398.666
$ws.Range("N68").Value = -18896.666

$ws.Range("H71").Value = 17398.666
$ws.Range("J71").Value = 17398.666
$ws.Range("L71").Value = 52195.99800000001
$ws.Range("N71").Value = -59683.99800000001

$ws.Range("H109").Value = 11360
$ws.Range("J109").Value = 11360
$ws.Range("L109").Value = 11360
$ws.Range("N109").Value = -13440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 458
$ws.Range("J12").Value = 641.7222
$ws.Range("L12").Value = 1925.1666
$ws.Range("N12").Value = -2271.1666

$ws.Range("H113").Value = 449993.4
$ws.Range("I113").Value = 448.42856
$ws.Range("J113").Value = 624816.4399999999
$ws.Range("K113").Value = 1345.28568
$ws.Range("L113").Value = 1874449.32
$ws.Range("M113").Value = 824.71432
$ws.Range("N113").Value = -1878789.32

$ws.Range("H132").Value = 428107.16
$ws.Range("I132").Value = 1317002.4
$ws.Range("J132").Value = 4823.7144
$ws.Range("K132").Value = 11853021.6
$ws.Range("L132").Value = 43413.4296
$ws.Range("M132").Value = -11850491.6
$ws.Range("N132").Value = -48473.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 34161.43
$ws.Range("J62").Value = 34161.43
$ws.Range("L62").Value = 34161.43
$ws.Range("N62").Value = -35533.43

$ws.Range("H63").Value = 23333.334
$ws.Range("J63").Value = 23333.334
$ws.Range("L63").Value = 23333.334
$ws.Range("N63").Value = -24705.334

$ws.Range("H65").Value = 34161.43
$ws.Range("J65").Value = 34161.43
$ws.Range("L65").Value = 102484.29
$ws.Range("N65").Value = -109348.29

$ws.Range("H66").Value = 23333.334
$ws.Range("J66").Value = 23333.334
$ws.Range("L66").Value = 70000.00199999999
$ws.Range("N66").Value = -76864.00199999999

$ws.Range("H97").Value = 1706.8462
$ws.Range("I97").Value = 1948.75
$ws.Range("J97").Value = 1319.8
$ws.Range("K97").Value = 1948.75
$ws.Range("L97").Value = 1319.8
$ws.Range("M97").Value = -1452.75
$ws.Range("N97").Value = -2311.8

$ws.Range("H123").Value = 35326
$ws.Range("J123").Value = 35326
$ws.Range("L123").Value = 35326
$ws.Range("N123").Value = -40226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1311.5555
$ws.Range("I22").Value = 1050.3334
$ws.Range("J22").Value = 1834
$ws.Range("K22").Value = 1050.3334
$ws.Range("L22").Value = 1834
$ws.Range("M22").Value = -755.3334
$ws.Range("N22").Value = -2424

$ws.Range("H27").Value = 1311.5555
$ws.Range("I27").Value = 1050.3334
$ws.Range("J27").Value = 1834
$ws.Range("K27").Value = 1050.3334
$ws.Range("L27").Value = 1834
$ws.Range("M27").Value = -943.3334
$ws.Range("N27").Value = -2048

$ws.Range("H132").Value = 3150.0264
$ws.Range("I132").Value = 3309.9333
$ws.Range("K132").Value = 9929.7999
$ws.Range("M132").Value = -7399.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1266.3889
$ws.Range("I96").Value = 1183.4445
$ws.Range("K96").Value = 1183.4445
$ws.Range("M96").Value = 189.5554999999999
